$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '43.710.77'
$ws.Range('E2').Value = '  +1.00%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.372.94'
$ws.Range('E3').Value = '  +6.15%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.13'
$ws.Range('E5').Value = '  +2.25%  '
$ws.Range('E6').Value = '  +2.55%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '73.62'
$ws.Range('E7').Value = '  +15.67%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.521'
$ws.Range('E9').Value = '  +19.36%  '
$ws.Range('E10').Value = '  +3.90%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '27.76'
$ws.Range('E11').Value = '  +4.62%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.724.09'
$ws.Range('E12').Value = '  +5.92%  '
$ws.Range('E13').Value = '  +2.42%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '16.40'
$ws.Range('E14').Value = '  +8.43%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.36'
$ws.Range('E15').Value = '  +6.30%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.882'
$ws.Range('E16').Value = '  +7.67%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.371.24'
$ws.Range('E17').Value = '  +6.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '43.573.62'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('E19').Value = '  +5.58%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '75.64'
$ws.Range('E20').Value = '  +3.92%  '
$ws.Range('E21').Value = '  +6.43%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '252.24'
$ws.Range('E22').Value = '  +2.89%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.80'
$ws.Range('E23').Value = '  +1.05%  '
$ws.Range('E24').Value = '  +0.03%  '
$ws.Range('E25').Value = '  +3.63%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '10.19'
$ws.Range('E26').Value = '  +5.64%  '
$ws.Range('E27').Value = '  -1.94%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '22.67'
$ws.Range('E28').Value = '  +5.17%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '172.89'
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('E30').Value = '  +9.75%  '
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  +3.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.08'
$ws.Range('E33').Value = '  +4.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0701'
$ws.Range('E34').Value = '  +4.41%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.10'
$ws.Range('E35').Value = '  +4.77%  '
$ws.Range('E36').Value = '  +6.39%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.69'
$ws.Range('E37').Value = '  +6.39%  '
$ws.Range('E38').Value = '  +9.53%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0259'
$ws.Range('E39').Value = '  +4.56%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '19.41'
$ws.Range('E40').Value = '  +15.07%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '8.96'
$ws.Range('E41').Value = '  +4.90%  '
$ws.Range('E42').Value = '  -0.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '100.59'
$ws.Range('E43').Value = '  +4.76%  '
$ws.Range('E44').Value = '  +11.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.51'
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('E46').Value = '  +4.13%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0960'
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.453.23'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('E49').Value = '  +9.53%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.597.55'
$ws.Range('E50').Value = '  +6.19%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.000204'
$ws.Range('E51').Value = '  -1.13%  '
